$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 11).Value = 2.12   # K6
$ws.Cells.Item(6, 16).Value = 2.95   # P6
$ws.Cells.Item(6, 34).Value = 15.5   # AH6
$ws.Cells.Item(6, 37).Value = 27   # AK6
$ws.Cells.Item(6, 39).Value = 450   # AM6
$ws.Cells.Item(6, 43).Value = 40   # AQ6
$ws.Cells.Item(6, 46).Value = 2.65   # AT6
$ws.Cells.Item(6, 53).Value = 110   # BA6
$ws.Cells.Item(7, 19).Value = 1.26   # S7
$ws.Cells.Item(7, 20).Value = 3.5   # T7
$ws.Cells.Item(15, 7).Value = 2.45   # G15
$ws.Cells.Item(15, 9).Value = 3.2   # I15
$ws.Cells.Item(15, 10).Value = 3.25   # J15
$ws.Cells.Item(15, 12).Value = 3.75   # L15
$ws.Cells.Item(15, 23).Value = 7   # W15
$ws.Cells.Item(15, 24).Value = 11   # X15
$ws.Cells.Item(15, 29).Value = 7   # AC15
$ws.Cells.Item(15, 30).Value = 5.5   # AD15
$ws.Cells.Item(15, 32).Value = 51   # AF15
$ws.Cells.Item(15, 35).Value = 12   # AI15
$ws.Cells.Item(15, 47).Value = 8.5   # AU15
$ws.Cells.Item(15, 49).Value = 4.75   # AW15
$ws.Cells.Item(15, 51).Value = 29   # AY15
$ws.Cells.Item(15, 54).Value = 251   # BB15
$ws.Cells.Item(16, 15).Value = 1.33   # O16
$ws.Cells.Item(16, 16).Value = 3.25   # P16
$ws.Cells.Item(17, 7).Value = 1.95   # G17
$ws.Cells.Item(17, 9).Value = 3.6   # I17
$ws.Cells.Item(17, 11).Value = 2.3   # K17
$ws.Cells.Item(17, 13).Value = 1.04   # M17
$ws.Cells.Item(17, 14).Value = 13   # N17
$ws.Cells.Item(17, 15).Value = 1.2   # O17
$ws.Cells.Item(17, 16).Value = 4.33   # P17
$ws.Cells.Item(17, 17).Value = 1.67   # Q17
$ws.Cells.Item(17, 18).Value = 2.15   # R17
$ws.Cells.Item(17, 23).Value = 9.5   # W17
$ws.Cells.Item(17, 26).Value = 17   # Z17
$ws.Cells.Item(17, 29).Value = 13   # AC17
$ws.Cells.Item(17, 35).Value = 13   # AI17
$ws.Cells.Item(17, 41).Value = 10   # AO17
$ws.Cells.Item(17, 42).Value = 19   # AP17
$ws.Cells.Item(17, 50).Value = 19   # AX17
$ws.Cells.Item(20, 7).Value = 2.4   # G20
$ws.Cells.Item(20, 8).Value = 2.88   # H20
$ws.Cells.Item(20, 9).Value = 3.3   # I20
$ws.Cells.Item(20, 10).Value = 3.1   # J20
$ws.Cells.Item(20, 11).Value = 2.05   # K20
$ws.Cells.Item(20, 12).Value = 3.75   # L20
$ws.Cells.Item(20, 13).Value = 1.08   # M20
$ws.Cells.Item(20, 14).Value = 7.5   # N20
$ws.Cells.Item(20, 15).Value = 1.36   # O20
$ws.Cells.Item(20, 16).Value = 3   # P20
$ws.Cells.Item(20, 17).Value = 2.25   # Q20
$ws.Cells.Item(20, 18).Value = 1.62   # R20
$ws.Cells.Item(20, 24).Value = 11   # X20
$ws.Cells.Item(20, 25).Value = 10   # Y20
$ws.Cells.Item(20, 26).Value = 23   # Z20
$ws.Cells.Item(20, 27).Value = 21   # AA20
$ws.Cells.Item(20, 29).Value = 7.5   # AC20
$ws.Cells.Item(20, 33).Value = 9   # AG20
$ws.Cells.Item(20, 36).Value = 34   # AJ20
$ws.Cells.Item(20, 37).Value = 29   # AK20
$ws.Cells.Item(20, 41).Value = 13   # AO20
$ws.Cells.Item(20, 44).Value = 67   # AR20
$ws.Cells.Item(20, 49).Value = 5   # AW20
$ws.Cells.Item(23, 7).Value = 2.63   # G23
$ws.Cells.Item(23, 8).Value = 3.75   # H23
$ws.Cells.Item(23, 10).Value = 3.1   # J23
$ws.Cells.Item(23, 17).Value = 1.53   # Q23
$ws.Cells.Item(23, 18).Value = 2.4   # R23
$ws.Cells.Item(23, 21).Value = 1.5   # U23
$ws.Cells.Item(23, 22).Value = 2.5   # V23
$ws.Cells.Item(23, 25).Value = 10   # Y23
$ws.Cells.Item(23, 28).Value = 21   # AB23
$ws.Cells.Item(23, 30).Value = 7.5   # AD23
$ws.Cells.Item(23, 32).Value = 34   # AF23
$ws.Cells.Item(23, 35).Value = 9.5   # AI23
$ws.Cells.Item(23, 50).Value = 13   # AX23
$ws.Cells.Item(23, 51).Value = 19   # AY23
$ws.Cells.Item(25, 7).Value = 2.15   # G25
$ws.Cells.Item(25, 8).Value = 3.45   # H25
$ws.Cells.Item(25, 10).Value = 2.72   # J25
$ws.Cells.Item(25, 11).Value = 2.22   # K25
$ws.Cells.Item(25, 15).Value = 1.24   # O25
$ws.Cells.Item(25, 16).Value = 3.7   # P25
$ws.Cells.Item(25, 17).Value = 1.72   # Q25
$ws.Cells.Item(25, 18).Value = 2.05   # R25
$ws.Cells.Item(25, 21).Value = 1.62   # U25
$ws.Cells.Item(25, 22).Value = 2.18   # V25
$ws.Cells.Item(25, 23).Value = 9   # W25
$ws.Cells.Item(25, 24).Value = 11.5   # X25
$ws.Cells.Item(25, 27).Value = 16.5   # AA25
$ws.Cells.Item(25, 28).Value = 24   # AB25
$ws.Cells.Item(25, 30).Value = 6.9   # AD25
$ws.Cells.Item(25, 32).Value = 50   # AF25
$ws.Cells.Item(25, 38).Value = 28   # AL25
$ws.Cells.Item(25, 41).Value = 11   # AO25
$ws.Cells.Item(25, 42).Value = 18   # AP25
$ws.Cells.Item(25, 48).Value = 55   # AV25
$ws.Cells.Item(25, 52).Value = 65   # AZ25
